$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05619466666666667
$ws.Range("I2").Value = 0.04986276087265156
$ws.Range("J2").Value = 0.07297477932340853
$ws.Range("M2").Value = 1.824475333333333
$ws.Range("N2").Value = 5.473426
$ws.Range("O2").Value = 0.1906606574278047
$ws.Range("P2").Value = 0.2015451970524477
$ws.Range("Q2").Value = 0.1025257831982222
$ws.Range("R2").Value = 0.9227320487840001
$ws.Range("S2").Value = 0.009506866769145165
$ws.Range("T2").Value = 0.01470771627859526

$ws.Range("G3").Value = 0.05619466666666667
$ws.Range("I3").Value = 0.04986276087265156
$ws.Range("J3").Value = 0.07297477932340853
$ws.Range("O3").Value = 0.6423822165107047
$ws.Range("P3").Value = 0.6790548829333741
$ws.Range("S3").Value = 0.03203095085071715
$ws.Range("T3").Value = 0.04955388023054599

$ws.Range("G4").Value = 0.05619466666666667
$ws.Range("I4").Value = 0.04986276087265156
$ws.Range("J4").Value = 0.07297477932340853
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009389666666666666
$ws.Range("N4").Value = 0.028169
$ws.Range("O4").Value = 0.0009812355294625031
$ws.Range("P4").Value = 0.001037252838673693
$ws.Range("Q4").Value = 0.0005276491884444444
$ws.Range("R4").Value = 0.004748842696
$ws.Range("S4").Value = 0.00004892711256533844
$ws.Range("T4").Value = 0.00007569329700479185

$ws.Range("G5").Value = 0.05619466666666667
$ws.Range("I5").Value = 0.04986276087265156
$ws.Range("J5").Value = 0.07297477932340853
$ws.Range("M5").Value = 1.5503715
$ws.Range("N5").Value = 3.100743
$ws.Range("O5").Value = 0.1620163583726162
$ws.Range("P5").Value = 0.1141770910840848
$ws.Range("Q5").Value = 0.08712260965200001
$ws.Range("R5").Value = 0.522735657912
$ws.Range("S5").Value = 0.008078582934991582
$ws.Range("T5").Value = 0.008332048025649803

$ws.Range("G6").Value = 0.05619466666666667
$ws.Range("I6").Value = 0.04986276087265156
$ws.Range("J6").Value = 0.07297477932340853
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03788966666666667
$ws.Range("N6").Value = 0.113669
$ws.Range("O6").Value = 0.003959532159411881
$ws.Range("P6").Value = 0.004185576091419648
$ws.Range("Q6").Value = 0.002129197188444445
$ws.Range("R6").Value = 0.019162774696
$ws.Range("S6").Value = 0.0001974332052323283
$ws.Range("T6").Value = 0.0003054414916126836

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.070792
$ws.Range("H7").Value = 2.141584
$ws.Range("I7").Value = 0.9501372391273485
$ws.Range("J7").Value = 0.9270252206765914
$ws.Range("M7").Value = 1.824475333333333
$ws.Range("N7").Value = 5.473426
$ws.Range("O7").Value = 0.1906606574278047
$ws.Range("P7").Value = 0.2015451970524477
$ws.Range("Q7").Value = 1.953633591130667
$ws.Range("R7").Value = 11.721801546784
$ws.Range("S7").Value = 0.1811537906586596
$ws.Range("T7").Value = 0.1868374807738524

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.070792
$ws.Range("H8").Value = 2.141584
$ws.Range("I8").Value = 0.9501372391273485
$ws.Range("J8").Value = 0.9270252206765914
$ws.Range("O8").Value = 0.6423822165107047
$ws.Range("P8").Value = 0.6790548829333741
$ws.Range("Q8").Value = 6.582267644784
$ws.Range("R8").Value = 39.493605868704
$ws.Range("S8").Value = 0.6103512656599875
$ws.Range("T8").Value = 0.6295010027028282

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.070792
$ws.Range("H9").Value = 2.141584
$ws.Range("I9").Value = 0.9501372391273485
$ws.Range("J9").Value = 0.9270252206765914
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.009389666666666666
$ws.Range("N9").Value = 0.028169
$ws.Range("O9").Value = 0.0009812355294625031
$ws.Range("P9").Value = 0.001037252838673693
$ws.Range("Q9").Value = 0.01005437994933333
$ws.Range("R9").Value = 0.060326279696
$ws.Range("S9").Value = 0.0009323084168971648
$ws.Range("T9").Value = 0.0009615595416689015

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.070792
$ws.Range("H10").Value = 2.141584
$ws.Range("I10").Value = 0.9501372391273485
$ws.Range("J10").Value = 0.9270252206765914
$ws.Range("M10").Value = 1.5503715
$ws.Range("N10").Value = 3.100743
$ws.Range("O10").Value = 0.1620163583726162
$ws.Range("P10").Value = 0.1141770910840848
$ws.Range("Q10").Value = 1.660125399228
$ws.Range("R10").Value = 6.640501596912
$ws.Range("S10").Value = 0.1539377754376246
$ws.Range("T10").Value = 0.105845043058435

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.070792
$ws.Range("H11").Value = 2.141584
$ws.Range("I11").Value = 0.9501372391273485
$ws.Range("J11").Value = 0.9270252206765914
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.03788966666666667
$ws.Range("N11").Value = 0.113669
$ws.Range("O11").Value = 0.003959532159411881
$ws.Range("P11").Value = 0.004185576091419648
$ws.Range("Q11").Value = 0.04057195194933334
$ws.Range("R11").Value = 0.243431711696
$ws.Range("S11").Value = 0.003762098954179553
$ws.Range("T11").Value = 0.003880134599806964
